# Rename the sheet from "Sheet" to "Favorite Things"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Favorite Things"

# Update column A: insert new header "Favorite Foods" and shift the
# existing food items (Pizza, Chocolate Cake, Broccoli) down one row.
$ws.Range("A1").Value = "Favorite Foods"
$ws.Range("A2").Value = "Pizza"
$ws.Range("A3").Value = "Chocolate Cake"
$ws.Range("A4").Value = "Broccoli"

# Add new column B with favorite colors.
$ws.Range("B1").Value = "Favorite Colors"
$ws.Range("B2").Value = "Blue"
$ws.Range("B3").Value = "Purple"
$ws.Range("B4").Value = "Green"
$ws.Range("B5").Value = "Orange"
